# Solution for "129. Sum Root to Leaf Numbers". zen-1
#
# The workbook is a LeetCode-progress tracker. Row 130 (problem #129,
# "Sum Root to Leaf Numbers") — and a handful of neighbouring rows whose
# "Finished" flag was stale — get flipped from "N" to "Y"; since the sheet
# has an active AutoFilter on Difficulty=Medium / Finished=N, those rows
# drop out of the visible set (hidden="1"). A few still-unsolved rows
# elsewhere get a yellow highlight. The active-cell selection also moves.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Finished" (column E) flips from N -> Y, and therefore get
# filtered out (hidden) because the AutoFilter only shows Difficulty=Medium
# AND Finished=N.
$doneRows = @(106, 122, 123, 130, 131, 132)

foreach ($r in $doneRows) {
    $ws.Range("E$r").Value = "Y"
    $ws.Rows.Item($r).Hidden = $true
}

# Highlight a few still-open rows in yellow.
$highlightCells = @("B157", "B162", "B164", "B168", "B177")
foreach ($addr in $highlightCells) {
    $ws.Range($addr).Interior.Color = 65535
}

# Move the active selection.
$ws.Range("B128").Select()
